# Insert a new row at position 13, shifting the existing rows 13-32 down to 14-33.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new price record.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 44729
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 100112043
$ws.Range("G13").Value = "Pepino dulce"
$ws.Range("H13").Value = "Cultivar IV Región"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 14000
$ws.Range("N13").Value = "$/bandeja 18 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 778
$ws.Range("Q13").Value = 18
$ws.Range("R13").Value = "Hortaliza"

# Match the date-formatted style used by the other rows in column D.
$ws.Range("D13").NumberFormat = $ws.Range("D14").NumberFormat
